# Weekly update: insert a new "Zanahoria" (Carrot) price record for
# "Vega Modelo de Temuco" ahead of the existing historical rows, shifting
# the previously-existing rows 354-393 down to 355-394.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 354; everything from the old row 354 onward
# (through the old last row, 393) shifts down by one row (to 355..394).
$ws.Rows(354).Insert()

# Populate the newly inserted row 354 with the new weekly data point.
$ws.Range("A354").Value = 10
$ws.Range("B354").Value = "Vega Modelo de Temuco"
$ws.Range("C354").Value = "La Araucanía"
$ws.Range("D354").Value = 44918
$ws.Range("E354").Value = 9
$ws.Range("F354").Value = 100114013
$ws.Range("G354").Value = "Zanahoria"
$ws.Range("H354").Value = "Sin especificar"
$ws.Range("I354").Value = "Primera"
$ws.Range("J354").Value = 250
$ws.Range("K354").Value = 7000
$ws.Range("L354").Value = 7000
$ws.Range("M354").Value = 7000
$ws.Range("N354").Value = "$/saco 20 kilos"
$ws.Range("O354").Value = "Provincia de Cautín"
$ws.Range("P354").Value = 350
$ws.Range("Q354").Value = 20
$ws.Range("R354").Value = "Hortaliza"
